$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.348
$ws.Range("AJ4").Value = 0.08599999999999999
$ws.Range("AK4").Value = 0.292
$ws.Range("AU4").Value = 0.196
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.166
$ws.Range("BA4").Value = 2.064
$ws.Range("BB4").Value = 0.143
$ws.Range("BC4").Value = 0.379
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.141
$ws.Range("BI4").Value = 0.375
$ws.Range("BM4").Value = 0.756
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.6879999999999999
$ws.Range("BQ4").Value = 0.765
$ws.Range("E4").Value = 0.423
$ws.Range("F4").Value = 0.07199999999999999
$ws.Range("G4").Value = 0.268
$ws.Range("N4").Value = 0.436
$ws.Range("O4").Value = 0.063
$ws.Range("P4").Value = 0.251
$ws.Range("Q4").Value = 0.22
$ws.Range("R4").Value = 0.111
$ws.Range("S4").Value = 0.333
$ws.Range("W4").Value = 0.345
$ws.Range("Y4").Value = 0.333
$ws.Range("AI5").Value = 0.372
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.308
$ws.Range("AU5").Value = 0.381
$ws.Range("AV5").Value = 0.091
$ws.Range("AW5").Value = 0.302
$ws.Range("BA5").Value = 1.3
$ws.Range("BB5").Value = 0.07199999999999999
$ws.Range("BC5").Value = 0.268
$ws.Range("BG5").Value = 0.382
$ws.Range("BH5").Value = 0.05
$ws.Range("BI5").Value = 0.223
$ws.Range("BM5").Value = 0.521
$ws.Range("BN5").Value = 0.047
$ws.Range("BO5").Value = 0.216
$ws.Range("BP5").Value = 0.433
$ws.Range("BQ5").Value = 0.456
$ws.Range("E5").Value = 0.538
$ws.Range("F5").Value = 0.08500000000000001
$ws.Range("G5").Value = 0.292
$ws.Range("N5").Value = 0.741
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.276
$ws.Range("Q5").Value = 0.145
$ws.Range("R5").Value = 0.045
$ws.Range("S5").Value = 0.212
$ws.Range("W5").Value = 0.333
$ws.Range("AI6").Value = 0.36
$ws.Range("AU6").Value = 0.259
$ws.Range("BA6").Value = 1.589
$ws.Range("BG6").Value = 0.501
$ws.Range("BM6").Value = 0.617
$ws.Range("BP6").Value = 0.53
$ws.Range("BQ6").Value = 0.569
$ws.Range("E6").Value = 0.474
$ws.Range("N6").Value = 0.549
$ws.Range("Q6").Value = 0.175
$ws.Range("W6").Value = 0.339
$ws.Range("AI7").Value = 0.367
$ws.Range("AU7").Value = 0.32
$ws.Range("BA7").Value = 1.402
$ws.Range("BG7").Value = 0.422
$ws.Range("BM7").Value = 0.556
$ws.Range("BP7").Value = 0.467
$ws.Range("BQ7").Value = 0.495
$ws.Range("E7").Value = 0.51
$ws.Range("N7").Value = 0.65
$ws.Range("Q7").Value = 0.156
$ws.Range("W7").Value = 0.335
$ws.Range("AI8").Value = 0.398
$ws.Range("AJ8").Value = 0.128
$ws.Range("AK8").Value = 0.358
$ws.Range("AU8").Value = 0.322
$ws.Range("AW8").Value = 0.29
$ws.Range("BA8").Value = 1.771
$ws.Range("BB8").Value = 0.108
$ws.Range("BC8").Value = 0.328
$ws.Range("BG8").Value = 0.5639999999999999
$ws.Range("BH8").Value = 0.108
$ws.Range("BI8").Value = 0.329
$ws.Range("BM8").Value = 0.675
$ws.Range("BN8").Value = 0.061
$ws.Range("BO8").Value = 0.247
$ws.Range("BP8").Value = 0.59
$ws.Range("BQ8").Value = 0.625
$ws.Range("E8").Value = 0.605
$ws.Range("F8").Value = 0.112
$ws.Range("G8").Value = 0.334
$ws.Range("N8").Value = 0.781
$ws.Range("O8").Value = 0.06
$ws.Range("P8").Value = 0.245
$ws.Range("Q8").Value = 0.159
$ws.Range("R8").Value = 0.08400000000000001
$ws.Range("S8").Value = 0.29
$ws.Range("W8").Value = 0.367
$ws.Range("X8").Value = 0.122
$ws.Range("Y8").Value = 0.35
$ws.Range("AI9").Value = 0.308
$ws.Range("AJ9").Value = 0.213
$ws.Range("AK9").Value = 0.462
$ws.Range("BA9").Value = 1.744
$ws.Range("BB9").Value = 0.25
$ws.Range("BC9").Value = 0.5
$ws.Range("BG9").Value = 0.605
$ws.Range("BH9").Value = 0.239
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.651
$ws.Range("BN9").Value = 0.227
$ws.Range("BO9").Value = 0.477
$ws.Range("BP9").Value = 0.581
$ws.Range("BQ9").Value = 0.619
$ws.Range("E9").Value = 0.553
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("N9").Value = 0.6879999999999999
$ws.Range("O9").Value = 0.215
$ws.Range("P9").Value = 0.463
$ws.Range("W9").Value = 0.247
$ws.Range("X9").Value = 0.186
$ws.Range("Y9").Value = 0.431
$ws.Range("AI10").Value = 0.436
$ws.Range("AJ10").Value = 0.246
$ws.Range("AK10").Value = 0.496
$ws.Range("AU10").Value = 0.311
$ws.Range("AV10").Value = 0.214
$ws.Range("AW10").Value = 0.463
$ws.Range("BA10").Value = 2.186
$ws.Range("BB10").Value = 0.211
$ws.Range("BC10").Value = 0.459
$ws.Range("BG10").Value = 0.674
$ws.Range("BH10").Value = 0.22
$ws.Range("BI10").Value = 0.469
$ws.Range("BM10").Value = 0.8139999999999999
$ws.Range("BN10").Value = 0.151
$ws.Range("BO10").Value = 0.389
$ws.Range("BP10").Value = 0.729
$ws.Range("BQ10").Value = 0.762
$ws.Range("E10").Value = 0.681
$ws.Range("F10").Value = 0.217
$ws.Range("G10").Value = 0.466
$ws.Range("N10").Value = 0.882
$ws.Range("O10").Value = 0.104
$ws.Range("P10").Value = 0.323
$ws.Range("W10").Value = 0.455
$ws.Range("X10").Value = 0.248
$ws.Range("Y10").Value = 0.498
$ws.Range("AI11").Value = 0.474
$ws.Range("AJ11").Value = 0.249
$ws.Range("AK11").Value = 0.499
$ws.Range("AU11").Value = 0.456
$ws.Range("AV11").Value = 0.248
$ws.Range("AW11").Value = 0.498
$ws.Range("BA11").Value = 2.186
$ws.Range("BB11").Value = 0.211
$ws.Range("BC11").Value = 0.459
$ws.Range("BG11").Value = 0.674
$ws.Range("BH11").Value = 0.22
$ws.Range("BI11").Value = 0.469
$ws.Range("BM11").Value = 0.8139999999999999
$ws.Range("BN11").Value = 0.151
$ws.Range("BO11").Value = 0.389
$ws.Range("BP11").Value = 0.729
$ws.Range("BQ11").Value = 0.768
$ws.Range("E11").Value = 0.713
$ws.Range("F11").Value = 0.205
$ws.Range("G11").Value = 0.452
$ws.Range("N11").Value = 0.903
$ws.Range("O11").Value = 0.08699999999999999
$ws.Range("P11").Value = 0.296
$ws.Range("W11").Value = 0.455
$ws.Range("X11").Value = 0.248
$ws.Range("Y11").Value = 0.498
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.767
$ws.Range("BB12").Value = 0.44
$ws.Range("BC12").Value = 0.663
$ws.Range("BG12").Value = 1.138
$ws.Range("BH12").Value = 0.188
$ws.Range("BI12").Value = 0.433
$ws.Range("BM12").Value = 1.229
$ws.Range("BN12").Value = 0.233
$ws.Range("BO12").Value = 0.483
$ws.Range("BP12").Value = 1.256
$ws.Range("BQ12").Value = 1.243
$ws.Range("E12").Value = 1.403
$ws.Range("F12").Value = 0.748
$ws.Range("G12").Value = 0.865
$ws.Range("N12").Value = 1.465
$ws.Range("O12").Value = 1.039
$ws.Range("P12").Value = 1.02
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.159
$ws.Range("BB13").Value = 0.277
$ws.Range("BC13").Value = 0.527
$ws.Range("BG13").Value = 0.542
$ws.Range("BH13").Value = 0.05
$ws.Range("BI13").Value = 0.224
$ws.Range("BM13").Value = 0.776
$ws.Range("BN13").Value = 0.164
$ws.Range("BO13").Value = 0.404
$ws.Range("BP13").Value = 0.72
$ws.Range("BQ13").Value = 0.661
$ws.Range("E13").Value = 1.573
$ws.Range("F13").Value = 0.652
$ws.Range("G13").Value = 0.8080000000000001
$ws.Range("N13").Value = 2.068
$ws.Range("O13").Value = 0.928
$ws.Range("P13").Value = 0.964
$ws.Range("W13").Value = 1.037
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
